# Updates the crypto price/volume table (and a couple of row reorderings)
# to match the latest scrape, per the "Updated cryptos list ... with
# GitHub Actions" commit.
#
# Numeric-looking price strings (e.g. "597.54") must stay TEXT, matching
# the source data's inlineStr cells - otherwise Excel auto-coerces them to
# floating point numbers (introducing binary rounding noise) and drops
# formatting like "1.00". We force text via NumberFormat "@" before the
# write, then reset the style back to Normal so no stray per-cell number
# format lingers afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.020.79'
$ws.Range('D3').Value = '3.308.44'
$ws.Range('E3').Value = '  +5.82%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.302.79'
$ws.Range('E8').Value = '  +5.75%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.471'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').Value = '3.850.92'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').Value = '3.306.85'
$ws.Range('E17').Value = '  +5.74%  '
$ws.Range('D18').Value = '64.069.12'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.743'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.14%  '
$ws.Range('E23').Value = '  +3.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.27%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.15%  '
$ws.Range('E31').Value = '  +2.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.63'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.61%  '
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.56'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('E35').Value = '  +1.45%  '
$ws.Range('E36').Value = '  +2.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.76%  '
$ws.Range('D38').Value = '0.0₃0736'
$ws.Range('E38').Value = '  +3.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0399'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '431.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.31%  '
$ws.Range('D41').Value = '3.017.85'
$ws.Range('E41').Value = '  +4.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('E44').Value = '  -6.19%  '
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.115'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +14.10%  '
